$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.199.86"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "3.129.52"
$ws.Range("E3").Value = "  -4.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "636.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.393"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.774"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "3.127.73"
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.557"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.178"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "88.993.39"
$ws.Range("E15").Value = "  -3.79%  "
$ws.Range("D16").Value = "3.704.78"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "32.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.37%  "
$ws.Range("D18").Value = "3.143.07"
$ws.Range("E18").Value = "  -3.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000225"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.78%  "
$ws.Range("E24").Value = "  -6.93%  "
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.14%  "
$ws.Range("D28").Value = "3.282.53"
$ws.Range("E28").Value = "  -4.83%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -10.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.977"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "504.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.01%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.147"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "21.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.18%  "
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.88%  "
$ws.Range("E44").Value = "  -7.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "144.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.131"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "163.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0648"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.44%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  -3.57%  "
